$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(4, 9).Value = "sv"
$ws.Cells.Item(4, 10).Value = "Statement-opinion"
$ws.Cells.Item(5, 9).Value = "ba"
$ws.Cells.Item(5, 10).Value = "Appreciation"
$ws.Cells.Item(8, 9).Value = "sv"
$ws.Cells.Item(8, 10).Value = "Statement-opinion"
$ws.Cells.Item(9, 9).Value = "ba"
$ws.Cells.Item(9, 10).Value = "Appreciation"
$ws.Cells.Item(11, 9).Value = "ba"
$ws.Cells.Item(11, 10).Value = "Appreciation"
$ws.Cells.Item(13, 9).Value = "ba"
$ws.Cells.Item(13, 10).Value = "Appreciation"
$ws.Cells.Item(16, 9).Value = "b"
$ws.Cells.Item(16, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(21, 9).Value = "aa"
$ws.Cells.Item(21, 10).Value = "Agree/Accept"
$ws.Cells.Item(24, 9).Value = "ba"
$ws.Cells.Item(24, 10).Value = "Appreciation"
$ws.Cells.Item(30, 9).Value = "ba"
$ws.Cells.Item(30, 10).Value = "Appreciation"
$ws.Cells.Item(32, 9).Value = "ba"
$ws.Cells.Item(32, 10).Value = "Appreciation"
$ws.Cells.Item(37, 9).Value = "ba"
$ws.Cells.Item(37, 10).Value = "Appreciation"
$ws.Cells.Item(44, 9).Value = "ba"
$ws.Cells.Item(44, 10).Value = "Appreciation"
$ws.Cells.Item(49, 9).Value = "b"
$ws.Cells.Item(49, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(52, 9).Value = "ba"
$ws.Cells.Item(52, 10).Value = "Appreciation"
$ws.Cells.Item(54, 9).Value = "aa"
$ws.Cells.Item(54, 10).Value = "Agree/Accept"
$ws.Cells.Item(58, 9).Value = "aa"
$ws.Cells.Item(58, 10).Value = "Agree/Accept"
$ws.Cells.Item(60, 9).Value = "sv"
$ws.Cells.Item(60, 10).Value = "Statement-opinion"
$ws.Cells.Item(63, 9).Value = "sd"
$ws.Cells.Item(63, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(82, 9).Value = "b"
$ws.Cells.Item(82, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(85, 9).Value = "ba"
$ws.Cells.Item(85, 10).Value = "Appreciation"
$ws.Cells.Item(90, 9).Value = "sv"
$ws.Cells.Item(90, 10).Value = "Statement-opinion"
$ws.Cells.Item(92, 9).Value = "sd"
$ws.Cells.Item(92, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(93, 9).Value = "ba"
$ws.Cells.Item(93, 10).Value = "Appreciation"
$ws.Cells.Item(94, 9).Value = "sd"
$ws.Cells.Item(94, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(98, 9).Value = "ba"
$ws.Cells.Item(98, 10).Value = "Appreciation"
$ws.Cells.Item(99, 9).Value = "sd"
$ws.Cells.Item(99, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(114, 9).Value = "%"
$ws.Cells.Item(114, 10).Value = "Uninterpretable"
$ws.Cells.Item(122, 9).Value = "b"
$ws.Cells.Item(122, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(126, 9).Value = "sd"
$ws.Cells.Item(126, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(129, 9).Value = "ba"
$ws.Cells.Item(129, 10).Value = "Appreciation"
$ws.Cells.Item(133, 9).Value = "ba"
$ws.Cells.Item(133, 10).Value = "Appreciation"
$ws.Cells.Item(135, 9).Value = "sd"
$ws.Cells.Item(135, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(137, 9).Value = "sd"
$ws.Cells.Item(137, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(143, 9).Value = "ba"
$ws.Cells.Item(143, 10).Value = "Appreciation"
